$wb = $excel.ActiveWorkbook

# Nature editor wants the figures renumbered: Fig1AB -> Fig2AB, Fig1C -> Fig2C
$wb.Worksheets.Item("Fig1AB").Name = "Fig2AB"
$wb.Worksheets.Item("Fig1C").Name = "Fig2C"

# Fig2C (formerly Fig1C) becomes the active tab, with a new cell selected
$wsC = $wb.Worksheets.Item("Fig2C")
$wsC.Activate()
$wsC.Range("D19").Select()
